$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new Price values are plain decimal numbers (e.g. "0.544"). The
# source data stores every Price/Volume cell as literal text, so entering
# such a value straight into a General-formatted cell would let Excel
# auto-convert it to a real number (and normalise away trailing zeros).
# Temporarily mark the cell as Text, assign the literal string, then
# restore the Normal style so the cell keeps its original (default)
# formatting while the stored value remains text.
$textCells = @("D5", "D6", "D7", "D8", "D9", "D11", "D15", "D16", "D18", "D22", "D23", "D25", "D26", "D27", "D28", "D30", "D35", "D37", "D38", "D42", "D43", "D46", "D48", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.369.05'
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").Value = '1.654.97'
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '213.28'
$ws.Range("E5").Value = '  -0.59%  '
$ws.Range("D6").Value = '0.544'
$ws.Range("E6").Value = '  +6.25%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("D8").Value = '23.53'
$ws.Range("E8").Value = '  +0.92%  '
$ws.Range("D9").Value = '0.262'
$ws.Range("E9").Value = '  +1.13%  '
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("D11").Value = '0.0905'
$ws.Range("E11").Value = '  +3.42%  '
$ws.Range("D12").Value = '1.889.66'
$ws.Range("E12").Value = '  -0.37%  '
$ws.Range("D13").Value = '1.651.62'
$ws.Range("E13").Value = '  -0.61%  '
$ws.Range("E14").Value = '  -1.07%  '
$ws.Range("D15").Value = '0.568'
$ws.Range("E15").Value = '  +3.84%  '
$ws.Range("D16").Value = '65.54'
$ws.Range("E16").Value = '  -0.43%  '
$ws.Range("D17").Value = '27.375.85'
$ws.Range("E17").Value = '  -0.57%  '
$ws.Range("D18").Value = '229.85'
$ws.Range("E18").Value = '  -6.99%  '
$ws.Range("D19").Value = '0.0₃0727'
$ws.Range("E19").Value = '  -0.61%  '
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("D22").Value = '4.35'
$ws.Range("E22").Value = '  -2.78%  '
$ws.Range("D23").Value = '9.27'
$ws.Range("E23").Value = '  -0.47%  '
$ws.Range("E24").Value = '  +0.49%  '
$ws.Range("D25").Value = '146.95'
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '0.115'
$ws.Range("E26").Value = '  +3.66%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '7.07'
$ws.Range("E27").Value = '  -1.39%  '
$ws.Range("D28").Value = '15.75'
$ws.Range("E28").Value = '  -2.72%  '
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("D30").Value = '0.0495'
$ws.Range("E30").Value = '  -0.70%  '
$ws.Range("E31").Value = '  -4.04%  '
$ws.Range("E32").Value = '  -0.87%  '
$ws.Range("E33").Value = '  +0.35%  '
$ws.Range("D34").Value = '1.422.03'
$ws.Range("E34").Value = '  -1.55%  '
$ws.Range("D35").Value = '1.56'
$ws.Range("E35").Value = '  +0.64%  '
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("D37").Value = '0.906'
$ws.Range("E37").Value = '  -2.45%  '
$ws.Range("D38").Value = '0.570'
$ws.Range("E38").Value = '  -1.78%  '
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("E40").Value = '  +0.79%  '
$ws.Range("E41").Value = '  -0.27%  '
$ws.Range("D42").Value = '5.55'
$ws.Range("E42").Value = '  +2.74%  '
$ws.Range("D43").Value = '65.08'
$ws.Range("E43").Value = '  -5.65%  '
$ws.Range("E44").Value = '  +0.67%  '
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = '1.798.16'
$ws.Range("E45").Value = '  -0.28%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").Value = '0.787'
$ws.Range("E46").Value = '  -0.30%  '
$ws.Range("E47").Value = '  -1.36%  '
$ws.Range("D48").Value = '88.12'
$ws.Range("E48").Value = '  -0.58%  '
$ws.Range("D49").Value = '0.0₆0105'
$ws.Range("E49").Value = '  -2.70%  '
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("D51").Value = '7.75'
$ws.Range("E51").Value = '  -0.74%  '

foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
